# Webex Dosen.xlsx - "Added all webex link"
# Fills in the remaining dosen (lecturer) names + their Webex meeting links,
# matching the order the original author entered them in (names for rows
# 3-8 first, then the links for rows 5,6,7,8,4,3 - with row 3's link made
# into a real Excel hyperlink - and finally name+link pairs for rows 9-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column A: names for rows 3-8 (entered first, as a block) ----
$ws.Range("A3").Value = "GUSMUL"
$ws.Range("A4").Value = "GUSDE"
$ws.Range("A5").Value = "MOGI"
$ws.Range("A6").Value = "GUNGDE"
$ws.Range("A7").Value = "HENDRA"
$ws.Range("A8").Value = "EKA"

# ---- Column B: links for rows 5,6,7,8,4 as plain text ----
$ws.Range("B5").Value = "https://universitas-udayana.webex.com/meet/arimogi"
$ws.Range("B6").Value = "https://universitas-udayana.webex.com/meet/gungde"
$ws.Range("B7").Value = "https://universitas-udayana.webex.com/meet/ibm.mahendra"
$ws.Range("B8").Value = "https://universitas-udayana.webex.com/meet/eka.karyawati"
$ws.Range("B4").Value = "https://universitas-udayana.webex.com/meet/dwidasmara"

# Row 3's link becomes a genuine hyperlink object (the only real one in
# the sheet) - leave the cell empty first so Add() seeds the text.
$ws.Hyperlinks.Add($ws.Range("B3"), "https://universitas-udayana.webex.com/meet/muliantara")

# ---- Remaining rows 9-14: name then link, one row at a time ----
$ws.Range("A9").Value = "ANOM"
$ws.Range("B9").Value = "https://universitas-udayana.webex.com/meet/anom.cp"

$ws.Range("A10").Value = "VIDA"
$ws.Range("B10").Value = "https://universitas-udayana.webex.com/meet/vida"

$ws.Range("A11").Value = "ASTUTI"
$ws.Range("B11").Value = "https://universitas-udayana.webex.com/meet/lg.astuti"

$ws.Range("A12").Value = "DAYU"
$ws.Range("B12").Value = "https://universitas-udayana.webex.com/meet/idgsuwiprabayantiputra"

$ws.Range("A13").Value = "COK"
$ws.Range("B13").Value = "https://universitas-udayana.webex.com/meet/cokorda"

$ws.Range("A14").Value = "SURYA"
$ws.Range("B14").Value = "https://universitas-udayana.webex.com/meet/surya"

# ---- Formatting: every filled-in link cell (B3:B14) gets a thin box
# border and the small Times New Roman font used throughout the sheet.
# (Row 3 keeps functioning as a real hyperlink underneath the formatting.)
$links = $ws.Range("B3:B14")
$links.Borders.LineStyle = 1
$links.Font.Name = "Times New Roman"
$links.Font.Size = 9
$links.Font.Underline = $false
$links.Font.ColorIndex = 1

# ---- Selection cursor ends on C17, matching the saved file ----
$ws.Range("C17").Select()
